# Apply Coastal CSPP (GA05MOAS-GL494) updates to the Omaha Cal Info workbook.

$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsAssetCal = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet -------------------------------------------------
# Reference designator corrected from GA05MOAS-GL001 to GA05MOAS-GL494
$wsMoorings.Range("A2").Value = "GA05MOAS-GL494"

# --- Asset_Cal_Info sheet --------------------------------------------
# Rows 2-5: FLORD reference designator
$wsAssetCal.Range("A2").Value = "GA05MOAS-GL494-01-FLORDM000"
$wsAssetCal.Range("A3").Value = "GA05MOAS-GL494-01-FLORDM000"
$wsAssetCal.Range("A4").Value = "GA05MOAS-GL494-01-FLORDM000"
$wsAssetCal.Range("A5").Value = "GA05MOAS-GL494-01-FLORDM000"

# Row 7: DOSTA reference designator
$wsAssetCal.Range("A7").Value = "GA05MOAS-GL494-02-DOSTAM000"

# Row 9: CTDGV reference designator
$wsAssetCal.Range("A9").Value = "GA05MOAS-GL494-04-CTDGVM000"

# Row 11: ENG reference designator
$wsAssetCal.Range("A11").Value = "GA05MOAS-GL494-00-ENG000000"

# --- Selected cell on Moorings sheet moved from D21 to B39 -----------
$wsMoorings.Activate()
$wsMoorings.Range("B39").Select()
